$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts Student Name .. RF ID Card No. one to the right)
$ws.Columns("B:B").Insert() | Out-Null

# New header for the inserted column
$ws.Range("B1").Value = "Admission No."

# Restore the width of the newly-inserted column to match its neighbour (was merged B:C before the insert)
$ws.Columns("B:B").ColumnWidth = 22

# Resize the last two columns (Alternative Mobile No. / RF ID Card No.) to fit their header text
$ws.Columns("H:H").ColumnWidth = 24.666666666666668
$ws.Columns("I:I").ColumnWidth = 16

# Move the active selection, matching the author's final cursor position
$ws.Range("G14").Select() | Out-Null
